$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 24 de Marzo de 2020 a las 21:16'

$ws.Range('B6').Value = 52431
$ws.Range('C6').Value = 8697
$ws.Range('E6').Value = 51387
$ws.Range('G6').Value = 121
$ws.Range('H6').Value = 674

$ws.Range('E8').Value = 29381
$ws.Range('G8').Value = 34
$ws.Range('H8').Value = 157

$ws.Range('B15').Value = 5283
$ws.Range('C15').Value = 809
$ws.Range('E15').Value = 5246

$ws.Range('B18').Value = 2590
$ws.Range('C18').Value = 499
$ws.Range('E18').Value = 2453

$ws.Range('B22').Value = 2100
$ws.Range('C22').Value = 176
$ws.Range('E22').Value = 2064

$ws.Range('A23').Value = 'Israel'
$ws.Range('B23').Value = 1930
$ws.Range('C23').Value = 488
$ws.Range('D23').Value = 53
$ws.Range('E23').Value = 1874
$ws.Range('F23').Value = 34
$ws.Range('G23').Value = 2
$ws.Range('H23').Value = 3

$ws.Range('A24').Value = 'Turquia'
$ws.Range('B24').Value = 1872
$ws.Range('C24').Value = 343
$ws.Range('D24').Value = 0
$ws.Range('E24').Value = 1828
$ws.Range('F24').Value = 0
$ws.Range('G24').Value = 7
$ws.Range('H24').Value = 44

$ws.Range('F50').Value = 9

$ws.Range('E69').Value = 206
$ws.Range('G69').Value = 1
$ws.Range('H69').Value = 2

$ws.Range('A77').Value = 'Bosnia y Herzegovina'
$ws.Range('B77').Value = 166
$ws.Range('C77').Value = 30
$ws.Range('D77').Value = 2
$ws.Range('E77').Value = 161
$ws.Range('F77').Value = 1
$ws.Range('G77').Value = 2
$ws.Range('H77').Value = 3

$ws.Range('A78').Value = 'Principado de Andorra'
$ws.Range('C78').Value = 31
$ws.Range('D78').Value = 1
$ws.Range('E78').Value = 162
$ws.Range('F78').Value = 7
$ws.Range('G78').Value = 0
$ws.Range('H78').Value = 1

$ws.Range('A119').Value = 'Ruanda'
$ws.Range('B119').Value = 40
$ws.Range('C119').Value = 4
$ws.Range('D119').Value = 0
$ws.Range('E119').Value = 40
$ws.Range('H119').Value = 0

$ws.Range('A120').Value = 'Puerto Rico'
$ws.Range('C120').Value = 8
$ws.Range('D120').Value = 1
$ws.Range('E120').Value = 36
$ws.Range('G120').Value = 0
$ws.Range('H120').Value = 2

$ws.Range('A121').Value = 'Banglades'
$ws.Range('B121').Value = 39
$ws.Range('C121').Value = 6
$ws.Range('D121').Value = 5
$ws.Range('E121').Value = 30
$ws.Range('G121').Value = 1
$ws.Range('H121').Value = 4

$ws.Range('A122').Value = 'Mayotte'
$ws.Range('C122').Value = 12

$ws.Range('A135').Value = 'Togo'
$ws.Range('C135').Value = 2
$ws.Range('D135').Value = 1
$ws.Range('H135').Value = 0

$ws.Range('A136').Value = 'Guyana'
$ws.Range('C136').Value = 0
$ws.Range('D136').Value = 0
$ws.Range('H136').Value = 1

$ws.Range('A137').Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range('C137').Value = 0

$ws.Range('A138').Value = 'Barbados'

$ws.Range('A139').Value = 'Madagascar'
$ws.Range('C139').Value = 5

$ws.Range('A150').Value = 'Seychelles'
$ws.Range('C150').Value = 0

$ws.Range('A151').Value = 'Surinam'
$ws.Range('C151').Value = 2

$ws.Range('A152').Value = 'Namibia'
$ws.Range('B152').Value = 7
$ws.Range('C152').Value = 3
$ws.Range('D152').Value = 2
$ws.Range('E152').Value = 5

$ws.Range('A155').Value = 'Benin'
$ws.Range('E155').Value = 6
$ws.Range('H155').Value = 0

$ws.Range('A156').Value = 'Gabon'
$ws.Range('C156').Value = 0

$ws.Range('A157').Value = 'Curazao'
$ws.Range('D157').Value = 0
$ws.Range('E157').Value = 5
$ws.Range('H157').Value = 1

$ws.Range('A158').Value = 'Islas Caimanes'
$ws.Range('B158').Value = 6
$ws.Range('C158').Value = 1
$ws.Range('H158').Value = 1

$ws.Range('A159').Value = 'El Salvador'
$ws.Range('C159').Value = 2
$ws.Range('D159').Value = 0
$ws.Range('E159').Value = 5

$ws.Range('A160').Value = 'Groenlandia'
$ws.Range('B160').Value = 5
$ws.Range('C160').Value = 1
$ws.Range('D160').Value = 2
$ws.Range('E160').Value = 3

$ws.Range('A161').Value = 'Congo'
$ws.Range('C161').Value = 0

$ws.Range('A162').Value = 'Guinea'

$ws.Range('A164').Value = 'Santa Sede'
$ws.Range('C164').Value = 3

$ws.Range('A166').Value = 'Fiyi'
$ws.Range('C166').Value = 1
$ws.Range('E166').Value = 4
$ws.Range('H166').Value = 0
